# Weekly update: insert a new price record for "Espinaca" (Femacal de La
# Calera) as row 41, pushing the existing rows 41-181 down to 42-182.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 41; Excel shifts every row
# at/after 41 down by one (old row 181 ends up at 182).
$ws.Rows(41).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44453
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 100112012
$ws.Range("G41").Value = "Espinaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 160
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = 3000
$ws.Range("N41").Value = "$/docena de atados (3 kilos)"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 1000
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = "Hortaliza"
